# Swap the "category" and "group" columns in the SectorGroup sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $d = $ws.Cells.Item($r, 4).Value()
    $e = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 4).Value = $e
    $ws.Cells.Item($r, 5).Value = $d

    $f = $ws.Cells.Item($r, 6).Value()
    $g = $ws.Cells.Item($r, 7).Value()
    $ws.Cells.Item($r, 6).Value = $g
    $ws.Cells.Item($r, 7).Value = $f
}
